$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (s=1, numFmt YYYY-MM-DD).
# Every populated row (2 through 96) had its date incremented by exactly one day:
# from serial 46075 (2026-02-22) to serial 46076 (2026-02-23).

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value2 = $v + 1
    }
}
